$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12
$ws.Range("H2").Value = 40

$ws.Range("G3").Value = 4
$ws.Range("H3").Value = 18

$ws.Range("G4").Value = 14
$ws.Range("H4").Value = 19

$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 3

$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 4

$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 10

$ws.Range("E9").Value = 26
$ws.Range("F9").Value = 11
$ws.Range("G9").Value = 9
$ws.Range("H9").Value = 20

$ws.Range("G10").Value = 2
$ws.Range("H10").Value = 3

$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 1

$ws.Range("G15").Value = 40
$ws.Range("H15").Value = 121

$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 8

$ws.Range("G17").Value = 32
$ws.Range("H17").Value = 86

$ws.Range("E18").Value = 108
$ws.Range("G18").Value = 36
$ws.Range("H18").Value = 81

$ws.Range("G19").Value = 13
$ws.Range("H19").Value = 40

$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 5

$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 2

$ws.Range("G24").Value = 4
$ws.Range("H24").Value = 17

$ws.Range("G25").Value = 8
$ws.Range("H25").Value = 19

$ws.Range("G26").Value = 10
$ws.Range("H26").Value = 21

$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 9

$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 14

$ws.Range("G29").Value = 3
$ws.Range("H29").Value = 12

$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 1

$ws.Range("G32").Value = 8
$ws.Range("H32").Value = 13

$ws.Range("G33").Value = 13
$ws.Range("H33").Value = 24

$ws.Range("G34").Value = 2
$ws.Range("H34").Value = 10

$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 3

$ws.Range("G36").Value = 32
$ws.Range("H36").Value = 74

$ws.Range("G37").Value = 12
$ws.Range("H37").Value = 40

$ws.Range("E38").Value = 72
$ws.Range("G38").Value = 20
$ws.Range("H38").Value = 35

$ws.Range("G39").Value = 8
$ws.Range("H39").Value = 22

$ws.Range("G40").Value = 3
$ws.Range("H40").Value = 13

$ws.Range("G41").Value = 11
$ws.Range("H41").Value = 27

$ws.Range("G42").Value = 9
$ws.Range("H42").Value = 26

$ws.Range("G43").Value = 3
$ws.Range("H43").Value = 17

$ws.Range("G44").Value = 10
$ws.Range("H44").Value = 24

$ws.Range("G45").Value = 7
$ws.Range("H45").Value = 21

$ws.Range("G46").Value = 9
$ws.Range("H46").Value = 18

$ws.Range("G47").Value = 11
$ws.Range("H47").Value = 46

$ws.Range("G48").Value = 5
$ws.Range("H48").Value = 24

$ws.Range("E49").Value = 69
$ws.Range("G49").Value = 17
$ws.Range("H49").Value = 56

$ws.Range("E50").Value = 24
$ws.Range("G50").Value = 8
$ws.Range("H50").Value = 16

$ws.Range("G53").Value = 2
$ws.Range("H53").Value = 4

$ws.Range("G55").Value = 3
$ws.Range("H55").Value = 6

$ws.Range("G56").Value = 3
$ws.Range("H56").Value = 5

$ws.Range("G57").Value = 4
$ws.Range("H57").Value = 7

$ws.Range("G59").Value = 4
$ws.Range("H59").Value = 8

$ws.Range("G60").Value = 5
$ws.Range("H60").Value = 14

$ws.Range("G61").Value = 10
$ws.Range("H61").Value = 20

$ws.Range("E62").Value = 40
$ws.Range("G62").Value = 14
$ws.Range("H62").Value = 24

$ws.Range("G63").Value = 8
$ws.Range("H63").Value = 20

$ws.Range("G64").Value = 5
$ws.Range("H64").Value = 22

$ws.Range("G65").Value = 13
$ws.Range("H65").Value = 22

$ws.Range("G66").Value = 8
$ws.Range("H66").Value = 29

$ws.Range("G67").Value = 9
$ws.Range("H67").Value = 30

$ws.Range("G68").Value = 4
$ws.Range("H68").Value = 13

$ws.Range("G69").Value = 3
$ws.Range("H69").Value = 11

$ws.Range("G70").Value = 12
$ws.Range("H70").Value = 32

$ws.Range("G71").Value = 10
$ws.Range("H71").Value = 26

$ws.Range("G72").Value = 11
$ws.Range("H72").Value = 31

$ws.Range("G73").Value = 12
$ws.Range("H73").Value = 23

$ws.Range("G74").Value = 4
$ws.Range("H74").Value = 11

$ws.Range("G75").Value = 5
$ws.Range("H75").Value = 12

$ws.Range("G76").Value = 17
$ws.Range("H76").Value = 34

$ws.Range("G77").Value = 17
$ws.Range("H77").Value = 36

$ws.Range("E78").Value = 45
$ws.Range("G78").Value = 21
$ws.Range("H78").Value = 40

$ws.Range("G79").Value = 12
$ws.Range("H79").Value = 30

$ws.Range("G80").Value = 12
$ws.Range("H80").Value = 21

$ws.Range("F81").Value = 7
$ws.Range("G81").Value = 5
$ws.Range("H81").Value = 12

$ws.Range("G82").Value = 6
$ws.Range("H82").Value = 8

$ws.Range("G83").Value = 7
$ws.Range("H83").Value = 9

$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 4

$ws.Range("G85").Value = 2
$ws.Range("H85").Value = 6

$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 2

$ws.Range("E87").Value = 15
$ws.Range("G87").Value = 7
$ws.Range("H87").Value = 10

$ws.Range("E88").Value = 20
$ws.Range("F88").Value = 11
$ws.Range("G88").Value = 8
$ws.Range("H88").Value = 19

$ws.Range("G89").Value = 7
$ws.Range("H89").Value = 21
